$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.712.26"
$ws.Range("E2").Value = "  +4.72%  "

# Row 3
$ws.Range("D3").Value = "2.250.48"
$ws.Range("E3").Value = "  +3.92%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "'248.93"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.13%  "

# Row 6
$ws.Range("D6").Value = "'0.633"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.03%  "

# Row 7
$ws.Range("D7").Value = "'70.57"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +7.28%  "

# Row 8
$ws.Range("E8").Value = "  -0.09%  "

# Row 9
$ws.Range("D9").Value = "'0.667"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +18.42%  "

# Row 10
$ws.Range("D10").Value = "'39.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +11.02%  "

# Row 11
$ws.Range("D11").Value = "'59.34"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.37%  "

# Row 12
$ws.Range("D12").Value = "'0.0965"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +4.91%  "

# Row 13
$ws.Range("D13").Value = "'7.51"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.78%  "

# Row 14
$ws.Range("E14").Value = "  +0.54%  "

# Row 15
$ws.Range("D15").Value = "2.580.01"
$ws.Range("E15").Value = "  +3.64%  "

# Row 16
$ws.Range("D16").Value = "'14.82"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.62%  "

# Row 17
$ws.Range("D17").Value = "'0.878"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.77%  "

# Row 18
$ws.Range("D18").Value = "2.249.70"
$ws.Range("E18").Value = "  +4.36%  "

# Row 19
$ws.Range("D19").Value = "42.660.93"
$ws.Range("E19").Value = "  +4.61%  "

# Row 20
$ws.Range("E20").Value = "  +5.77%  "

# Row 21
$ws.Range("E21").Value = "  +3.52%  "

# Row 22
$ws.Range("D22").Value = "'72.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.30%  "

# Row 23
$ws.Range("D23").Value = "'235.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.13%  "

# Row 24
$ws.Range("E24").Value = "  -0.55%  "

# Row 25
$ws.Range("D25").Value = "'3.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.53%  "

# Row 26
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'11.57"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.22%  "

# Row 27
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D27").Value = "'1.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "

# Row 28
$ws.Range("E28").Value = "  +1.31%  "

# Row 29
$ws.Range("E29").Value = "  -1.47%  "

# Row 30
$ws.Range("E30").Value = "  -0.73%  "

# Row 31
$ws.Range("D31").Value = "'167.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.28%  "

# Row 32
$ws.Range("D32").Value = "'20.86"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.82%  "

# Row 33
$ws.Range("D33").Value = "'6.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +15.14%  "

# Row 34
$ws.Range("E34").Value = "  +6.81%  "

# Row 35
$ws.Range("E35").Value = "  +7.92%  "

# Row 36
$ws.Range("D36").Value = "'30.99"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +22.28%  "

# Row 37
$ws.Range("E37").Value = "  +4.24%  "

# Row 38
$ws.Range("D38").Value = "'4.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.21%  "

# Row 39
$ws.Range("D39").Value = "'4.70"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.94%  "

# Row 40
$ws.Range("E40").Value = "  +7.88%  "

# Row 41
$ws.Range("D41").Value = "'2.30"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.45%  "

# Row 42
$ws.Range("D42").Value = "'12.50"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +8.22%  "

# Row 43
$ws.Range("D43").Value = "'5.80"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.60%  "

# Row 44
$ws.Range("D44").Value = "'62.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.28%  "

# Row 45
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").Value = "'0.201"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +6.78%  "

# Row 46
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "'9.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.83%  "

# Row 47
$ws.Range("E47").Value = "  +1.30%  "

# Row 48
$ws.Range("E48").Value = "  +3.53%  "

# Row 49
$ws.Range("E49").Value = "  -0.70%  "

# Row 50
$ws.Range("E50").Value = "  +1.05%  "

# Row 51
$ws.Range("E51").Value = "  +4.55%  "
